$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update input cells per the "hw - led board redesign" change.
# All other changed cells in the diff (B9, C9, D9, D10, B11, C11, D11, C13)
# are formula-driven and recalculate automatically from these inputs.
$ws.Range("B1").Value = 0
$ws.Range("B2").Value = 0.35
$ws.Range("B3").Value = 9
$ws.Range("C10").Value = 0.8

$excel.Calculate()
